# Gestion des livres completed
# Populate the first row of the active sheet with the book record:
# id, author, title, quantity, price, isbn (all stored as text values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting only on the numeric-looking fields so they are
# kept as strings (ids / quantity / price / isbn are not real numbers),
# matching the source data.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("D1").NumberFormat = "@"
$ws.Range("E1").NumberFormat = "@"
$ws.Range("F1").NumberFormat = "@"

$ws.Range("A1").Value = "874596"
$ws.Range("B1").Value = "monpierre"
$ws.Range("C1").Value = "davenci"
$ws.Range("D1").Value = "5"
$ws.Range("E1").Value = "100"
$ws.Range("F1").Value = "1253468025641"
